$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# V2G row (row 4): update Trading rev. (C4), Adj. rev. (F4) and dSoH (ppm)
# (G4). These numeric-looking values are stored as text in the sheet, so
# briefly switch to a text number format while assigning (otherwise Excel
# auto-converts the literal to a real number), then clear the formatting
# back off again so the cell ends up with the same (default) style as
# before the edit - only its text content changes.
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "-1"
$ws.Range("C4").ClearFormats()

$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "-1.6"
$ws.Range("F4").ClearFormats()

$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "61.2"
$ws.Range("G4").ClearFormats()
